$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.844.52"
$ws.Range("E2").Value = "  -1.12%  "

# Row 3
$ws.Range("D3").Value = "1.856.25"
$ws.Range("E3").Value = "  -0.53%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'304.34"
$ws.Range("E5").Value = "  -0.83%  "

# Row 6
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.08%  "

# Row 7
$ws.Range("D7").Value = "'0.5042"
$ws.Range("E7").Value = "  -2.29%  "

# Row 8
$ws.Range("E8").Value = "  -2.66%  "

# Row 9
$ws.Range("D9").Value = "'0.07170"
$ws.Range("E9").Value = "  +0.12%  "

# Row 10
$ws.Range("D10").Value = "'0.8910"
$ws.Range("E10").Value = "  +0.71%  "

# Row 11
$ws.Range("D11").Value = "'20.63"
$ws.Range("E11").Value = "  -0.19%  "

# Row 12
$ws.Range("D12").Value = "'0.07517"
$ws.Range("E12").Value = "  -0.52%  "

# Row 13
$ws.Range("D13").Value = "1.840.88"
$ws.Range("E13").Value = "  -1.41%  "

# Row 14
$ws.Range("D14").Value = "'92.05"
$ws.Range("E14").Value = "  +3.09%  "

# Row 15
$ws.Range("D15").Value = "'5.226"
$ws.Range("E15").Value = "  -1.93%  "

# Row 16
$ws.Range("E16").Value = "  +0.06%  "

# Row 17
$ws.Range("D17").Value = "'0.000008499"

# Row 18
$ws.Range("E18").Value = "  -0.65%  "

# Row 19
$ws.Range("E19").Value = "  +0.02%  "

# Row 20
$ws.Range("D20").Value = "26.871.05"
$ws.Range("E20").Value = "  -1.20%  "

# Row 21
$ws.Range("D21").Value = "'5.027"
$ws.Range("E21").Value = "  +0.00%  "

# Row 22
$ws.Range("D22").Value = "2.089.41"
$ws.Range("E22").Value = "  -1.35%  "

# Row 23
$ws.Range("D23").Value = "'10.33"
$ws.Range("E23").Value = "  -2.54%  "

# Row 24
$ws.Range("D24").Value = "'6.454"
$ws.Range("E24").Value = "  -0.27%  "

# Row 25
$ws.Range("D25").Value = "'146.46"
$ws.Range("E25").Value = "  -2.95%  "

# Row 26
$ws.Range("D26").Value = "'1.799"
$ws.Range("E26").Value = "  -2.58%  "

# Row 27
$ws.Range("E27").Value = "  -1.06%  "

# Row 28
$ws.Range("D28").Value = "'2.057"
$ws.Range("E28").Value = "  -4.38%  "

# Row 29
$ws.Range("D29").Value = "'112.91"

# Row 30
$ws.Range("D30").Value = "'4.636"
$ws.Range("E30").Value = "  -2.24%  "

# Row 31
$ws.Range("D31").Value = "'4.658"
$ws.Range("E31").Value = "  -0.63%  "

# Row 32
$ws.Range("E32").Value = "  +1.98%  "

# Row 33
$ws.Range("D33").Value = "'0.05087"
$ws.Range("E33").Value = "  -1.45%  "

# Row 34
$ws.Range("D34").Value = "'2.991"
$ws.Range("E34").Value = "  -3.45%  "

# Row 35
$ws.Range("D35").Value = "'0.7361"
$ws.Range("E35").Value = "  -2.24%  "

# Row 36
$ws.Range("D36").Value = "'1.145"
$ws.Range("E36").Value = "  -2.33%  "

# Row 37
$ws.Range("D37").Value = "'3.236"
$ws.Range("E37").Value = "  +7.01%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'2.513"
$ws.Range("E38").Value = "  -0.56%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01992"
$ws.Range("E39").Value = "  -1.99%  "

# Row 40
$ws.Range("D40").Value = "'1.073"
$ws.Range("E40").Value = "  -0.71%  "

# Row 41
$ws.Range("D41").Value = "'0.5328"
$ws.Range("E41").Value = "  -0.44%  "

# Row 42
$ws.Range("D42").Value = "'119.20"
$ws.Range("E42").Value = "  +3.46%  "

# Row 43
$ws.Range("D43").Value = "'6.477"
$ws.Range("E43").Value = "  -2.64%  "

# Row 44
$ws.Range("D44").Value = "'8.348"

# Row 45
$ws.Range("E45").Value = "  -1.16%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.4636"
$ws.Range("E46").Value = "  -0.60%  "

# Row 47
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "  +0.10%  "

# Row 48
$ws.Range("D48").Value = "'9.929"
$ws.Range("E48").Value = "  -1.94%  "

# Row 49
$ws.Range("D49").Value = "'1.556"
$ws.Range("E49").Value = "  -1.08%  "

# Row 50
$ws.Range("D50").Value = "'36.88"
$ws.Range("E50").Value = "  +1.46%  "

# Row 51
$ws.Range("D51").Value = "'62.77"
$ws.Range("E51").Value = "  -3.28%  "
